$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "['0.96 (0.92, 0.99)']"
$ws.Range("E2").Value = "['0.67 (0.59, 0.76)']"
$ws.Range("D3").Value = "['0.89 (0.79, 0.99)']"
$ws.Range("E3").Value = "['0.42 (0.3, 0.71)']"
